# MicroDustHeroConfig.xlsx: rename the hero-data worksheet so the config
# (and anything that builds from it) refers to it as "Heros" instead of
# the old placeholder name "UnitProto".
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "UnitProto") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    # Fallback: original sheet already renamed or not found by that name -
    # just use the first (only) sheet in the workbook.
    $ws = $wb.Worksheets.Item(1)
}

$ws.Name = "Heros"
